$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.3205216666666666
$ws.Range("H2").Value = 0.961565
$ws.Range("I2").Value = 0.1355444593128396
$ws.Range("J2").Value = 0.1355444593128396
$ws.Range("M2").Value = 10.823698
$ws.Range("N2").Value = 32.471094
$ws.Range("O2").Value = 0.3079767696785641
$ws.Range("P2").Value = 0.3079767696785641
$ws.Range("Q2").Value = 3.469229722456666
$ws.Range("R2").Value = 31.22306750211
$ws.Range("S2").Value = 0.0417445447269959
$ws.Range("T2").Value = 0.04174454472699589

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.3205216666666666
$ws.Range("H3").Value = 0.961565
$ws.Range("I3").Value = 0.1355444593128396
$ws.Range("J3").Value = 0.1355444593128396
$ws.Range("O3").Value = 0.5341523443640567
$ws.Range("P3").Value = 0.5341523443640566
$ws.Range("Q3").Value = 6.017003137352777
$ws.Range("R3").Value = 54.153028236175
$ws.Range("S3").Value = 0.07240139070751175
$ws.Range("T3").Value = 0.07240139070751174

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.3205216666666666
$ws.Range("H4").Value = 0.961565
$ws.Range("I4").Value = 0.1355444593128396
$ws.Range("J4").Value = 0.1355444593128396
$ws.Range("M4").Value = 5.548297666666667
$ws.Range("N4").Value = 16.644893
$ws.Range("O4").Value = 0.1578708859573793
$ws.Range("P4").Value = 0.1578708859573793
$ws.Range("Q4").Value = 1.778349615282778
$ws.Range("R4").Value = 16.005146537545
$ws.Range("S4").Value = 0.02139852387833194
$ws.Range("T4").Value = 0.02139852387833193

$ws.Range("I5").Value = 0.4633580360449179
$ws.Range("J5").Value = 0.4633580360449179
$ws.Range("M5").Value = 10.823698
$ws.Range("N5").Value = 32.471094
$ws.Range("O5").Value = 0.3079767696785641
$ws.Range("P5").Value = 0.3079767696785641
$ws.Range("Q5").Value = 11.85954393809667
$ws.Range("R5").Value = 106.73589544287
$ws.Range("S5").Value = 0.1427035111457175
$ws.Range("T5").Value = 0.1427035111457175

$ws.Range("I6").Value = 0.4633580360449179
$ws.Range("J6").Value = 0.4633580360449179
$ws.Range("O6").Value = 0.5341523443640567
$ws.Range("P6").Value = 0.5341523443640566
$ws.Range("S6").Value = 0.247503781233318
$ws.Range("T6").Value = 0.2475037812333179

$ws.Range("I7").Value = 0.4633580360449179
$ws.Range("J7").Value = 0.4633580360449179
$ws.Range("M7").Value = 5.548297666666667
$ws.Range("N7").Value = 16.644893
$ws.Range("O7").Value = 0.1578708859573793
$ws.Range("P7").Value = 0.1578708859573793
$ws.Range("Q7").Value = 6.079279000529444
$ws.Range("R7").Value = 54.713511004765
$ws.Range("S7").Value = 0.07315074366588251
$ws.Range("T7").Value = 0.07315074366588249

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.9484743333333334
$ws.Range("H8").Value = 2.845423
$ws.Range("I8").Value = 0.4010975046422426
$ws.Range("J8").Value = 0.4010975046422425
$ws.Range("M8").Value = 10.823698
$ws.Range("N8").Value = 32.471094
$ws.Range("O8").Value = 0.3079767696785641
$ws.Range("P8").Value = 0.3079767696785641
$ws.Range("Q8").Value = 10.26599974475134
$ws.Range("R8").Value = 92.39399770276201
$ws.Range("S8").Value = 0.1235287138058508
$ws.Range("T8").Value = 0.1235287138058507

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.9484743333333334
$ws.Range("H9").Value = 2.845423
$ws.Range("I9").Value = 0.4010975046422426
$ws.Range("J9").Value = 0.4010975046422425
$ws.Range("O9").Value = 0.5341523443640567
$ws.Range("P9").Value = 0.5341523443640566
$ws.Range("Q9").Value = 17.80526445752056
$ws.Range("R9").Value = 160.247380117685
$ws.Range("S9").Value = 0.214247172423227
$ws.Range("T9").Value = 0.2142471724232269

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.9484743333333334
$ws.Range("H10").Value = 2.845423
$ws.Range("I10").Value = 0.4010975046422426
$ws.Range("J10").Value = 0.4010975046422425
$ws.Range("M10").Value = 5.548297666666667
$ws.Range("N10").Value = 16.644893
$ws.Range("O10").Value = 0.1578708859573793
$ws.Range("P10").Value = 0.1578708859573793
$ws.Range("Q10").Value = 5.262417930526556
$ws.Range("R10").Value = 47.361761374739
$ws.Range("S10").Value = 0.0633216184131649
$ws.Range("T10").Value = 0.06332161841316489
